$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.039.59'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.865.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.25'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5085'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3746'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07149'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8851'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.65'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07584'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.858.71'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.311'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.41'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9995'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008420'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.08'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9989'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.075.19'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.033'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.092.08'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.51'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.470'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.84'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.96'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.104'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.75'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.674'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.717'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09062'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05144'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.049'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.155'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7287'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02042'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.94%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.481'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.46%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.038'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5309'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.549'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.12'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.280'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9990'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4627'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.982'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.566'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.57'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.89'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.95%  '